$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Sheet6 (2)" is the active/selected tab in this workbook

# --- Row 7: replace the old "C1 / NA" row with the new SLOT summary row ---
$ws.Range("A7").Value = "Slot1L: Slot3R;Slot2L: Slot1R;Slot3L: Slot4R;Slot4L: Slot 2R"
$ws.Range("B7").Value = "NA"
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = " "
$ws.Range("E7").Value = 1

# --- Row 8 ---
$ws.Range("A8").Value = "Slot1L " + [char]0x2260 + " Slot3R"
$ws.Range("B8").Value = "FD1"
$ws.Range("C8").Value = " "
$ws.Range("D8").Value = " "
$ws.Range("E8").Value = 2

# --- Row 9 ---
$ws.Range("A9").Value = "Slot2L " + [char]0x2260 + " Slot1R"
$ws.Range("B9").Value = "FD2"
$ws.Range("C9").Value = " "
$ws.Range("D9").Value = " "
$ws.Range("E9").Value = 3

# --- Row 10 ---
$ws.Range("A10").Value = "Slot3L " + [char]0x2260 + " Slot4R"
$ws.Range("B10").Value = "FD3"
$ws.Range("C10").Value = " "
$ws.Range("D10").Value = " "
$ws.Range("E10").Value = 4

# --- Row 11 (new row) ---
$ws.Range("A11").Value = "Slot4L " + [char]0x2260 + " Slot2R"
$ws.Range("B11").Value = "FD4"
$ws.Range("C11").Value = " "
$ws.Range("D11").Value = " "
$ws.Range("E11").Value = 5

# Wrap text (claims cellXfs style index 2) + row heights for the rebuilt table rows
$ws.Range("A7:B11").WrapText = $true
$ws.Rows.Item(7).RowHeight = 80
$ws.Rows.Item(8).RowHeight = 16
$ws.Rows.Item(9).RowHeight = 16
$ws.Rows.Item(10).RowHeight = 16
$ws.Rows.Item(11).RowHeight = 16

# Explicit "General" number format on the sequence column (claims cellXfs style index 3)
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "General"

# --- Update the summary block (changed last, matching the order new text was introduced) ---
$ws.Range("C2").Value = "SLOT"
$ws.Range("B4").Value = "SLOT Test"
